# Scheduled runner update: refresh market price / profit columns (H-N)
# across all item sheets, per latest Market Board snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 10600.167
$ws.Range("I18").Value = 10150.25
$ws.Range("K18").Value = 10150.25
$ws.Range("M18").Value = -9866.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 313.0476
$ws.Range("I33").Value = 214.42105
$ws.Range("J33").Value = 1250
$ws.Range("K33").Value = 214.42105
$ws.Range("L33").Value = 1250
$ws.Range("M33").Value = 14.57894999999999
$ws.Range("N33").Value = -1708

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 978.1667
$ws.Range("I103").Value = 773.8
$ws.Range("J103").Value = 2000
$ws.Range("K103").Value = 2321.4
$ws.Range("L103").Value = 6000
$ws.Range("M103").Value = -1735.4
$ws.Range("N103").Value = -7172

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 26250.514
$ws.Range("I116").Value = 30347.691
$ws.Range("K116").Value = 30347.691
$ws.Range("M116").Value = -26905.691

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 533551.5600000001
$ws.Range("I137").Value = 914.2593000000001
$ws.Range("J137").Value = 912004.4
$ws.Range("K137").Value = 2742.7779
$ws.Range("L137").Value = 2736013.2
$ws.Range("M137").Value = -192.7779
$ws.Range("N137").Value = -2741113.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2391.7576
$ws.Range("I138").Value = 1574.921
$ws.Range("K138").Value = 4724.763
$ws.Range("M138").Value = 415.2370000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5005084
$ws.Range("I32").Value = 5379897
$ws.Range("J32").Value = 25424.285
$ws.Range("K32").Value = 5379897
$ws.Range("L32").Value = 25424.285
$ws.Range("M32").Value = -5379610
$ws.Range("N32").Value = -25998.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6562.1113
$ws.Range("I45").Value = 3671.3845
$ws.Range("K45").Value = 3671.3845
$ws.Range("M45").Value = -3294.3845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H49").Value = 21499.9
$ws.Range("J49").Value = 21499.9
$ws.Range("L49").Value = 21499.9
$ws.Range("N49").Value = -22019.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 22500
$ws.Range("J54").Value = 22500
$ws.Range("L54").Value = 22500
$ws.Range("N54").Value = -24038

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2559.7188
$ws.Range("I102").Value = 2074.64
$ws.Range("K102").Value = 2074.64
$ws.Range("M102").Value = -452.6399999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3037.35
$ws.Range("I122").Value = 1459.091
$ws.Range("K122").Value = 4377.272999999999
$ws.Range("M122").Value = -1927.272999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1788.9524
$ws.Range("J86").Value = 2023.4546
$ws.Range("L86").Value = 2023.4546
$ws.Range("N86").Value = -4269.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1788.9524
$ws.Range("J89").Value = 2023.4546
$ws.Range("L89").Value = 10117.273
$ws.Range("N89").Value = -21349.273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2315.842
$ws.Range("I99").Value = 1190.25
$ws.Range("J99").Value = 4245.4287
$ws.Range("K99").Value = 1190.25
$ws.Range("L99").Value = 4245.4287
$ws.Range("M99").Value = 307.75
$ws.Range("N99").Value = -7241.4287

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1425.7727
$ws.Range("I107").Value = 1473.2565
$ws.Range("K107").Value = 1473.2565
$ws.Range("M107").Value = 446.7435

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2227.4285
$ws.Range("I86").Value = 2253.889
$ws.Range("J86").Value = 2179.8
$ws.Range("K86").Value = 2253.889
$ws.Range("L86").Value = 2179.8
$ws.Range("M86").Value = -1130.889
$ws.Range("N86").Value = -4425.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 2227.4285
$ws.Range("I89").Value = 2253.889
$ws.Range("J89").Value = 2179.8
$ws.Range("K89").Value = 11269.445
$ws.Range("L89").Value = 10899
$ws.Range("M89").Value = -5653.445
$ws.Range("N89").Value = -22131

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 518.9259
$ws.Range("I107").Value = 518.4400000000001
$ws.Range("K107").Value = 518.4400000000001
$ws.Range("M107").Value = 1401.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 55833370
$ws.Range("I2").Value = 416695.9
$ws.Range("K2").Value = 2500175.4
$ws.Range("M2").Value = -2500062.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 20.571428
$ws.Range("J12").Value = 26.333334
$ws.Range("L12").Value = 79.00000199999999
$ws.Range("N12").Value = -425.000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 262.9
$ws.Range("I92").Value = 267.625
$ws.Range("J92").Value = 244
$ws.Range("K92").Value = 802.875
$ws.Range("L92").Value = 732
$ws.Range("M92").Value = 445.125
$ws.Range("N92").Value = -3228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 310.33334
$ws.Range("I97").Value = 87.333336
$ws.Range("J97").Value = 533.3333
$ws.Range("K97").Value = 262.000008
$ws.Range("L97").Value = 1599.9999
$ws.Range("M97").Value = 233.999992
$ws.Range("N97").Value = -2591.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3240.1428
$ws.Range("I109").Value = 1697.4615
$ws.Range("J109").Value = 5747
$ws.Range("K109").Value = 5092.3845
$ws.Range("L109").Value = 17241
$ws.Range("M109").Value = -4052.3845
$ws.Range("N109").Value = -19321

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 4998.5
$ws.Range("I141").Value = 4998.5
$ws.Range("K141").Value = 14995.5
$ws.Range("M141").Value = -9815.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1624.7838
$ws.Range("I97").Value = 282.61905
$ws.Range("K97").Value = 282.61905
$ws.Range("M97").Value = 213.38095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 360.76923
$ws.Range("I107").Value = 270.22223
$ws.Range("K107").Value = 270.22223
$ws.Range("M107").Value = 1649.77777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 36492.38
$ws.Range("I22").Value = 111973.445
$ws.Range("J22").Value = 2525.9
$ws.Range("K22").Value = 111973.445
$ws.Range("L22").Value = 2525.9
$ws.Range("M22").Value = -111678.445
$ws.Range("N22").Value = -3115.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 36492.38
$ws.Range("I27").Value = 111973.445
$ws.Range("J27").Value = 2525.9
$ws.Range("K27").Value = 111973.445
$ws.Range("L27").Value = 2525.9
$ws.Range("M27").Value = -111866.445
$ws.Range("N27").Value = -2739.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3217.923
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 3277.75
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 3277.75
$ws.Range("M46").Value = -2312
$ws.Range("N46").Value = -3653.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 718.75
$ws.Range("I55").Value = 190.25
$ws.Range("J55").Value = 1511.5
$ws.Range("K55").Value = 190.25
$ws.Range("L55").Value = 1511.5
$ws.Range("M55").Value = -17.25
$ws.Range("N55").Value = -1857.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 2740
$ws.Range("I23").Value = 600
$ws.Range("J23").Value = 4166.6665
$ws.Range("K23").Value = 600
$ws.Range("L23").Value = 4166.6665
$ws.Range("M23").Value = -371
$ws.Range("N23").Value = -4624.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 551449.5
$ws.Range("I132").Value = 825402
$ws.Range("J132").Value = 3544.5715
$ws.Range("K132").Value = 2476206
$ws.Range("L132").Value = 10633.7145
$ws.Range("M132").Value = -2473676
$ws.Range("N132").Value = -15693.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 11181171
$ws.Range("I136").Value = 14620614
$ws.Range("K136").Value = 43861842
$ws.Range("M136").Value = -43859292
